# "added search combos to CapStatement"
#
# The "sps" worksheet has a search-parameter "combos" column (Q) which
# listed single modifiers per search parameter ("any"/"none"/etc.). This
# change turns it into a "combo_pairs" column describing actual
# multi-parameter search combinations (comma separated parameter lists),
# removing the entries that no longer apply and updating the remaining
# ones with the new combo text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sps")

# Header: "combos" -> "combo_pairs" (R1 "combo_conf" is unchanged)
$ws.Range("Q1").Value = "combo_pairs"

# Questionnaire._id / Questionnaire.url no longer have a combo entry
$ws.Range("Q2").ClearContents()
$ws.Range("Q3").ClearContents()

# Questionnaire.publisher combo
$ws.Range("Q6").Value = "context-type,version"

# QuestionnaireResponse._id / QuestionnaireResponse.questionnaire no longer
# have a combo entry
$ws.Range("Q9").ClearContents()
$ws.Range("Q10").ClearContents()

# QuestionnaireResponse.patient combo (text unchanged, rewritten so the
# shared-string table reflects the new content correctly)
$ws.Range("Q11").Value = "author,source,context"

# QuestionnaireResponse.status combo
$ws.Range("Q13").Value = "questionnaire,patient,context,status,author,source"

# Questionnaire.status combo
$ws.Range("Q4").Value = "title,publisher,version,context-type-value"

# Widen column Q so the new (longer) combo text is fully visible, matching
# the author's manual column resize when they reviewed the new data.
$ws.Columns.Item(17).ColumnWidth = 44.6

# Reflect where the author was looking at on each sheet when they saved.
# "resources" must stay the active tab (as in the original file), so it is
# selected last.
$ws.Range("Q7").Select() | Out-Null

$wsOps = $wb.Worksheets.Item("ops")
$wsOps.Range("A1:C1").Select() | Out-Null

$wsResources = $wb.Worksheets.Item("resources")
$wsResources.Range("D2").Select() | Out-Null
